$wb = $excel.ActiveWorkbook

function Add-MarketSheet($SourceName, $AfterName, $NewName, $TicketValue, $MarketValue) {
    $source = $wb.Worksheets.Item($SourceName)
    $after = $wb.Worksheets.Item($AfterName)
    $source.Copy($null, $after)

    $newSheet = $wb.Worksheets.Item($after.Index + 1)
    $newSheet.Name = $NewName
    # Match the shared-string insertion order observed in the target workbook:
    # the ticket reference is written before the market name on each new sheet.
    $newSheet.Range("B4").Value = $TicketValue
    $newSheet.Range("B2").Value = $MarketValue
}

# Add the three new market sheets after Austria, cloning the Denmark sheet as
# the closest existing template (same layout/merges/styles).
Add-MarketSheet "Denmark" "Austria" "Russia" "NGC-2929/T2901" "Russia Market"
Add-MarketSheet "Denmark" "Russia" "Finland" "NGC-3130/T2944" "Finland Market"
Add-MarketSheet "Denmark" "Finland" "Hungary" "NGC-3104/T2993" "Hungary Market"

# Move the active tab/selection from Austria to Croatia, matching the
# workbook-level activeTab change (19 -> 15) seen in the target edit.
$croatia = $wb.Worksheets.Item("Croatia")
$croatia.Activate()
$croatia.Range("E20").Select()
